# Adds a new "RULE 10" validation entry to both the Summary and Details
# sheets of the CLFS sample input validation report.

$wb = $excel.ActiveWorkbook

# --- Summary sheet -------------------------------------------------------
# Insert a new row above the existing row 5 ("RULE 5") so the new "RULE 10"
# entry appears first, and the old "RULE 5" row shifts down to row 6.
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Rows.Item(5).Insert()

$wsSummary.Cells.Item(5, 1).Value = "RULE 10"
$wsSummary.Cells.Item(5, 2).Value = "Was your main job last week a paid internship, traineeship or apprenticeship? & Type of Employment?"
$wsSummary.Cells.Item(5, 3).Value = "Internship/Traineeship/Apprenticeship must be Fixed-Term contract employee"
$wsSummary.Cells.Item(5, 4).Value = 1

# --- Details sheet ---------------------------------------------------------
# Append a new row 8 describing the same RULE 10 violation found for a
# specific respondent.
$wsDetails = $wb.Worksheets.Item("Details")

$wsDetails.Cells.Item(8, 1).Value = "CLFS_sample_input.xlsx"
$wsDetails.Cells.Item(8, 2).Value = 4
$wsDetails.Cells.Item(8, 3).Value = "697c2c580deae81fbb49c180"
$wsDetails.Cells.Item(8, 4).Value = 1
$wsDetails.Cells.Item(8, 5).Value = "Chen Jia Hui"
$wsDetails.Cells.Item(8, 6).Value = "RULE 10"
$wsDetails.Cells.Item(8, 7).Value = "Was your main job last week a paid internship, traineeship or apprenticeship? & Type of Employment?"
$wsDetails.Cells.Item(8, 8).Value = "Internship/Traineeship/Apprenticeship must be Fixed-Term contract employee"
